$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New test-case rows 8-10 (TC_07/08/09 "Check ... DD in New Ticket page") ---

# Row 8: finish filling in the (previously blank) row 8 with TC_07 data.
# Column A needs the "bordered" data-row style already used by A2/A3/A4/A5/A7.
$ws.Range("A2").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A8").Value = "TC_07_Check Priority DD in New Ticket page"
$ws.Range("B8").Value = "admin"
$ws.Range("C8").Value = "admin"
$ws.Range("D8").Value = "Aqua"

# Row 9: brand-new row, same style pattern as row 8 (A col bordered data style,
# B:I plain bordered style already used across the sheet).
$ws.Range("A2").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("B2:I2").Copy()
$ws.Range("B9:I9").PasteSpecial(-4122)
$ws.Range("A9").Value = "TC_08_Check Severity DD in New Ticket page"
$ws.Range("B9").Value = "admin"
$ws.Range("C9").Value = "admin"
$ws.Range("D9").Value = "orange"

# Row 10: brand-new row; column A re-uses the border-less style from A6.
$ws.Range("A6").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("B2:I2").Copy()
$ws.Range("B10:I10").PasteSpecial(-4122)
$ws.Range("A10").Value = "TC_09_Check category DD in New Ticket page"
$ws.Range("B10").Value = "admin"
$ws.Range("C10").Value = "admin"
$ws.Range("D10").Value = "nature"

# --- Sheet cosmetics that came along with the edit ---

# Column A got a bit wider to fit the longer scenario names.
$ws.Columns.Item(1).ColumnWidth = 47.498697916666664

# Selection/scroll moved as the author clicked around after adding the rows.
$ws.Range("F14").Select()

$excel.CutCopyMode = $false
